$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Y")

# Insert two new columns before column D, shifting existing data (D:K) to (F:M)
$ws.Range("D:E").EntireColumn.Insert()

# Copy number formats from column F (old D, now shifted right) back onto new D and E columns
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("E5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarters of data (columns D and E)
$ws.Range("D7").Value = "12/31/2018"
$ws.Range("E7").Value = "09/30/2018"
$ws.Range("D8").Value = 1969700
$ws.Range("E8").Value = 1807200
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 800
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 7300
$ws.Range("E15").Value = 5500
$ws.Range("D17").Value = 2093900
$ws.Range("E17").Value = 1805400
$ws.Range("D18").Value = -124200
$ws.Range("E18").Value = 1800
$ws.Range("D20").Value = -741800
$ws.Range("E20").Value = 370200
$ws.Range("D21").Value = -834100
$ws.Range("E21").Value = 403000
$ws.Range("D22").Value = 24700
$ws.Range("E22").Value = 22200
$ws.Range("D23").Value = -890700
$ws.Range("E23").Value = 349900
$ws.Range("D24").Value = -186300
$ws.Range("E24").Value = 60400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -704400
$ws.Range("E26").Value = 289500
$ws.Range("D27").Value = -712100
$ws.Range("E27").Value = 284900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 741800
$ws.Range("E32").Value = -370200
$ws.Range("D33").Value = -712100
$ws.Range("E33").Value = 284900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -712100
$ws.Range("E35").Value = 284900
$ws.Range("D38").Value = "12/31/2018"
$ws.Range("E38").Value = "09/30/2018"
$ws.Range("D41").Value = 1038800
$ws.Range("E41").Value = 646900
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 842600
$ws.Range("E43").Value = 839000
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 17639700
$ws.Range("E47").Value = 19142900
$ws.Range("D48").Value = 195200
$ws.Range("E48").Value = 196300
$ws.Range("D49").Value = 1008300
$ws.Range("E49").Value = 811900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 164900
$ws.Range("E52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 25344900
$ws.Range("E54").Value = 25795900
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 14686000
$ws.Range("E59").Value = 14316300
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 1669000
$ws.Range("E61").Value = 1581700
$ws.Range("D62").Value = "NA"
$ws.Range("E62").Value = 6200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 17652200
$ws.Range("E66").Value = 17200800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 5577400
$ws.Range("E72").Value = 6289400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 7692700
$ws.Range("E76").Value = 8595100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = "12/31/2018"
$ws.Range("E80").Value = "09/30/2018"
$ws.Range("D81").Value = -712100
$ws.Range("E81").Value = 284900
$ws.Range("D83").Value = 31900
$ws.Range("E83").Value = 30900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 78600
$ws.Range("E89").Value = 174600
$ws.Range("D91").Value = 15600
$ws.Range("E91").Value = -25700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 476000
$ws.Range("E94").Value = -115400
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -159400
$ws.Range("E100").Value = -31000
$ws.Range("D101").Value = -3400
$ws.Range("E101").Value = -4600
$ws.Range("D102").Value = 391900
$ws.Range("E102").Value = 23700
$ws.Range("F62").Value = "NA"
$ws.Range("G62").Value = "NA"
$ws.Range("H62").Value = "NA"
$ws.Range("I62").Value = "NA"
$ws.Range("J62").Value = "NA"
